$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A4").NumberFormat = "@"

$ws.Range("A2").Value = "0010003000008       "
$ws.Range("B2").Value = "PASTILLAS FRENOS ACERO MDA01 STL. HPD02-PAR       "
$ws.Range("C2").Value = 7

$ws.Range("A3").Value = "0010003000018       "
$ws.Range("B3").Value = "FRENOS VBRAKE ALUMINIO TWISTER HVK140 NEGRO - SET "
$ws.Range("C3").Value = 201

$ws.Range("A4").Value = "0010003000065       "
$ws.Range("B4").Value = "ZAPATAS CRS STL RA01 PLUS                         "
$ws.Range("C4").Value = 3
